# Update crypto price/volume snapshot (GitHub Actions scheduled refresh).
# D column = Price, E column = Volume(1h) change. Values that look numeric
# are prefixed with a leading apostrophe so Excel keeps them as literal text
# (matching the original inline-string cells) instead of coercing to Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.702.99"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "2.613.19"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'574.82"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("D6").Value = "'156.42"
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -3.49%  "
$ws.Range("E9").Value = "  -7.01%  "
$ws.Range("D10").Value = "'5.85"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").Value = "'0.381"
$ws.Range("E11").Value = "  -5.26%  "
$ws.Range("D12").Value = "'0.156"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'28.27"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "3.077.60"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("E15").Value = "  -8.80%  "
$ws.Range("D16").Value = "63.570.39"
$ws.Range("E16").Value = "  -3.38%  "
$ws.Range("D17").Value = "2.612.84"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("E18").Value = "  -4.84%  "
$ws.Range("D19").Value = "'7.60"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("E20").Value = "  -5.84%  "
$ws.Range("D21").Value = "'343.55"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'67.57"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").Value = "'1.81"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  -4.15%  "
$ws.Range("D26").Value = "'9.17"
$ws.Range("E26").Value = "  -5.25%  "
$ws.Range("D27").Value = "'585.21"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "'1.59"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "'0.161"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("E33").Value = "  -4.49%  "
$ws.Range("D34").Value = "'6.59"
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").Value = "'5.36"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("D37").Value = "'19.76"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'154.37"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("D42").Value = "'2.53"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'41.43"
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").Value = "'157.38"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").Value = "'23.76"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "'3.90"
$ws.Range("E46").Value = "  -5.28%  "
$ws.Range("D47").Value = "'0.0590"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").Value = "'0.0247"
$ws.Range("E50").Value = "  -5.05%  "
$ws.Range("D51").Value = "'18.90"
$ws.Range("E51").Value = "  -4.41%  "
